$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (it currently sits between
#    "color " and "tints" in the "Changing layout theme..." bullet).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2. Turn the last bullet ("Tasklist appearance when not loggedin") into
#    "Task title disappearing".
$rng = $d.Content
$rng.Find.Execute("Tasklist appearance when not loggedin", $true, $false, $false, $false, $false, $true, 1, $false, "Task title disappearing", 2) | Out-Null

# Re-touch the same text so any leftover proofing-error markers around the
# old runs get cleared up (self replace is a no-op textually).
$rng2 = $d.Content
$rng2.Find.Execute("Task title disappearing", $true, $false, $false, $false, $false, $true, 1, $false, "Task title disappearing", 2) | Out-Null

# 3. Add a new bullet after it, in the same list, with the second issue.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "When checked as done, list changing back to something else"

# 4. Re-create the "_GoBack" bookmark at the very end of the new bullet.
#    Inserting a bookmark exactly at end-of-paragraph is unreliable, so we
#    temporarily extend the paragraph with a placeholder character, anchor
#    the bookmark right before it, then remove the placeholder again.
$endPara = $d.Paragraphs.Last
$tailRange = $endPara.Range
$tailRange.Collapse(0)
$tailRange.InsertAfter("X")

$bmRange = $d.Content
$bmRange.Find.Execute("elseX", $true, $false, $false, $false, $false, $true, 1, $false, "else", 2) | Out-Null

$anchorRange = $d.Content
$anchorRange.Find.Execute("something else", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $anchorRange) | Out-Null
